$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$words = @("It", "also", "works", "if", "you", "want", "to", "mount", "a", "range", "to", "an", "entire", "column.")

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 5).Value = $words[$i]
}

$ws.Activate()
$ws.Range("F9").Select()
